# Apply scheduled-runner updates to Sheets/Durandal_Profits.xlsx (Leve profit tables).
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific rows
# across the ALC, ARM, BSM, CRP, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3591110.8
$ws.Range("I98").Value = 7742.3335
$ws.Range("K98").Value = 7742.3335
$ws.Range("M98").Value = -6244.3335
$ws.Range("H122").Value = 3591110.8
$ws.Range("I122").Value = 7742.3335
$ws.Range("K122").Value = 23227.0005
$ws.Range("M122").Value = -20777.0005
$ws.Range("H133").Value = 45327.145
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 45327.145
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 45327.145
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -55447.145
$ws.Range("H137").Value = 1206.8572
$ws.Range("I137").Value = 942
$ws.Range("J137").Value = 1637.25
$ws.Range("K137").Value = 2826
$ws.Range("L137").Value = 4911.75
$ws.Range("M137").Value = -276
$ws.Range("N137").Value = -10011.75
$ws.Range("H140").Value = 68761.42999999999
$ws.Range("J140").Value = 88266
$ws.Range("L140").Value = 88266
$ws.Range("N140").Value = -98626

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 466878.88
$ws.Range("I32").Value = 5831.559
$ws.Range("J32").Value = 9534143
$ws.Range("K32").Value = 5831.559
$ws.Range("L32").Value = 9534143
$ws.Range("M32").Value = -5544.559
$ws.Range("N32").Value = -9534717
$ws.Range("H74").Value = 1105
$ws.Range("I74").Value = 660.4
$ws.Range("K74").Value = 660.4
$ws.Range("M74").Value = 213.6
$ws.Range("H77").Value = 1105
$ws.Range("I77").Value = 660.4
$ws.Range("K77").Value = 3302
$ws.Range("M77").Value = 1066
$ws.Range("H88").Value = 5187.375
$ws.Range("I88").Value = 1499.6666
$ws.Range("K88").Value = 1499.6666
$ws.Range("M88").Value = -1093.6666
$ws.Range("H91").Value = 5187.375
$ws.Range("I91").Value = 1499.6666
$ws.Range("K91").Value = 1499.6666
$ws.Range("M91").Value = -95.66660000000002
$ws.Range("H106").Value = 69800
$ws.Range("J106").Value = 69800
$ws.Range("L106").Value = 69800
$ws.Range("N106").Value = -72324
$ws.Range("H132").Value = 1913.9811
$ws.Range("I132").Value = 1567.7174
$ws.Range("J132").Value = 4189.4287
$ws.Range("K132").Value = 4703.1522
$ws.Range("L132").Value = 12568.2861
$ws.Range("M132").Value = -2173.1522
$ws.Range("N132").Value = -17628.2861
$ws.Range("H133").Value = 77532.625
$ws.Range("J133").Value = 77532.625
$ws.Range("L133").Value = 77532.625
$ws.Range("N133").Value = -82592.625
$ws.Range("H140").Value = 93119.086
$ws.Range("J140").Value = 93119.086
$ws.Range("L140").Value = 93119.086
$ws.Range("N140").Value = -103479.086
$ws.Range("H141").Value = 59800
$ws.Range("J141").Value = 62707.69
$ws.Range("L141").Value = 62707.69
$ws.Range("N141").Value = -73067.69

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H140").Value = 88900
$ws.Range("J140").Value = 88900
$ws.Range("L140").Value = 88900
$ws.Range("N140").Value = -99260

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4274.614
$ws.Range("I31").Value = 2683.0256
$ws.Range("J31").Value = 5541.3877
$ws.Range("K31").Value = 2683.0256
$ws.Range("L31").Value = 5541.3877
$ws.Range("M31").Value = -2388.0256
$ws.Range("N31").Value = -6131.3877
$ws.Range("H34").Value = 4274.614
$ws.Range("I34").Value = 2683.0256
$ws.Range("J34").Value = 5541.3877
$ws.Range("K34").Value = 2683.0256
$ws.Range("L34").Value = 5541.3877
$ws.Range("M34").Value = -2481.0256
$ws.Range("N34").Value = -5945.3877
$ws.Range("H36").Value = 3893
$ws.Range("I36").Value = 3563.4285
$ws.Range("K36").Value = 3563.4285
$ws.Range("M36").Value = -3175.4285
$ws.Range("H40").Value = 3893
$ws.Range("I40").Value = 3563.4285
$ws.Range("K40").Value = 3563.4285
$ws.Range("M40").Value = -3403.4285
$ws.Range("H57").Value = 8161
$ws.Range("J57").Value = 8161
$ws.Range("L57").Value = 8161
$ws.Range("N57").Value = -9281
$ws.Range("H58").Value = 2262.3333
$ws.Range("I58").Value = 600
$ws.Range("K58").Value = 600
$ws.Range("M58").Value = -397
$ws.Range("H132").Value = 2405
$ws.Range("I132").Value = 1501.7368
$ws.Range("J132").Value = 4856.7144
$ws.Range("K132").Value = 4505.2104
$ws.Range("L132").Value = 14570.1432
$ws.Range("M132").Value = -1975.2104
$ws.Range("N132").Value = -19630.1432
$ws.Range("H136").Value = 2262.3333
$ws.Range("I136").Value = 600
$ws.Range("K136").Value = 1800
$ws.Range("M136").Value = 750
$ws.Range("H140").Value = 89950
$ws.Range("J140").Value = 89950
$ws.Range("L140").Value = 89950
$ws.Range("N140").Value = -100310

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 1999.4445
$ws.Range("J24").Value = 1999.4445
$ws.Range("L24").Value = 1999.4445
$ws.Range("N24").Value = -2345.4445
$ws.Range("H132").Value = 3101.476
$ws.Range("I132").Value = 2055
$ws.Range("K132").Value = 6165
$ws.Range("M132").Value = -3635
$ws.Range("H133").Value = 54800
$ws.Range("J133").Value = 54800
$ws.Range("L133").Value = 54800
$ws.Range("N133").Value = -64920
$ws.Range("H138").Value = 68300
$ws.Range("J138").Value = 68300
$ws.Range("L138").Value = 68300
$ws.Range("N138").Value = -78580

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2520.7896
$ws.Range("I40").Value = 2289.875
$ws.Range("J40").Value = 3752.3333
$ws.Range("K40").Value = 2289.875
$ws.Range("L40").Value = 3752.3333
$ws.Range("M40").Value = -2153.875
$ws.Range("N40").Value = -4024.3333
$ws.Range("H57").Value = 303000
$ws.Range("I57").Value = 401666.66
$ws.Range("K57").Value = 401666.66
$ws.Range("M57").Value = -401100.66
$ws.Range("H122").Value = 3121.9
$ws.Range("I122").Value = 2733.6428
$ws.Range("J122").Value = 4027.8333
$ws.Range("K122").Value = 8200.928400000001
$ws.Range("L122").Value = 12083.4999
$ws.Range("M122").Value = -5750.928400000001
$ws.Range("N122").Value = -16983.4999
$ws.Range("H132").Value = 4947.2
$ws.Range("I132").Value = 8700.5
$ws.Range("J132").Value = 3180.9412
$ws.Range("K132").Value = 26101.5
$ws.Range("L132").Value = 9542.8236
$ws.Range("M132").Value = -23571.5
$ws.Range("N132").Value = -14602.8236
$ws.Range("H138").Value = 48420.43
$ws.Range("J138").Value = 48420.43
$ws.Range("L138").Value = 48420.43
$ws.Range("N138").Value = -58700.43

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 3249.1667
$ws.Range("I37").Value = 1500
$ws.Range("K37").Value = 1500
$ws.Range("M37").Value = -1297
$ws.Range("H122").Value = 2987
$ws.Range("I122").Value = 2987
$ws.Range("K122").Value = 8961
$ws.Range("M122").Value = -6511
$ws.Range("H132").Value = 25426794
$ws.Range("I132").Value = 31915918
$ws.Range("K132").Value = 95747754
$ws.Range("M132").Value = -95745224
$ws.Range("H133").Value = 21729
$ws.Range("J133").Value = 21729
$ws.Range("L133").Value = 21729
$ws.Range("N133").Value = -31849
$ws.Range("H136").Value = 822.3226
$ws.Range("I136").Value = 765.9655
$ws.Range("J136").Value = 1639.5
$ws.Range("K136").Value = 2297.8965
$ws.Range("L136").Value = 4918.5
$ws.Range("M136").Value = 252.1035000000002
$ws.Range("N136").Value = -10018.5
